# [Abraham]: fixed select all in enterprise.
#
# The "select all" / export routine in the Enterprise flow concatenated the
# IFI codes for a row into a single comma separated cell (e.g.
# "1002003004001,1002003004002"). That plain string was ambiguous for the
# consumer, so it is now wrapped in square brackets to look like the JSON
# array the importer expects: "[1002003004001,1002003004002]".
#
# This updates the two sample/demo rows on the "Clientes" sheet (W2 and W3)
# that hard-code that value, turns off the wrap-text formatting that had
# been forcing the long code onto two lines, and leaves the selection on
# the last edited cell (W3), matching what Excel records after manually
# editing that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "[1002003004001,1002003004002]"

$w2 = $ws.Range("W2")
$w2.Value = $newValue
$w2.WrapText = $false

$w3 = $ws.Range("W3")
$w3.Value = $newValue
$w3.WrapText = $false

# Leave the selection where the author last clicked after the edit.
$null = $w3.Select()
